# Scheduled-runner data refresh for the Jenova market-board leve-profit
# workbook: pulls fresh currentAveragePrice* figures per crafting job
# table and recomputes the dependent Leve price/profit columns.
#
# iron_native COM-interop script (PowerShell-style). The workbook is
# already open as $excel.ActiveWorkbook; each job table lives on its
# own worksheet tab (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 143.53847
$ws.Range("I12").Value = 145.57143
$ws.Range("K12").Value = 145.57143
$ws.Range("M12").Value = 24.42857000000001
# Row 88
$ws.Range("H88").Value = 4747.143
$ws.Range("I88").Value = 4867.3335
$ws.Range("J88").Value = 4657
$ws.Range("K88").Value = 4867.3335
$ws.Range("L88").Value = 4657
$ws.Range("M88").Value = -4461.3335
$ws.Range("N88").Value = -5469
# Row 91
$ws.Range("H91").Value = 4747.143
$ws.Range("I91").Value = 4867.3335
$ws.Range("J91").Value = 4657
$ws.Range("K91").Value = 4867.3335
$ws.Range("L91").Value = 4657
$ws.Range("M91").Value = -3463.3335
$ws.Range("N91").Value = -7465
# Row 137
$ws.Range("H137").Value = 16614.5
$ws.Range("I137").Value = 1669.625
$ws.Range("K137").Value = 5008.875
$ws.Range("M137").Value = -2458.875
# Row 138
$ws.Range("H138").Value = 4840.6787
$ws.Range("J138").Value = 6473.543
$ws.Range("L138").Value = 19420.629
$ws.Range("N138").Value = -29700.629

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1467.7368
$ws.Range("I32").Value = 1534.5769
$ws.Range("J32").Value = 772.6
$ws.Range("K32").Value = 1534.5769
$ws.Range("L32").Value = 772.6
$ws.Range("M32").Value = -1247.5769
$ws.Range("N32").Value = -1346.6
# Row 132
$ws.Range("H132").Value = 230450.55
$ws.Range("I132").Value = 374406.7
$ws.Range("K132").Value = 1123220.1
$ws.Range("M132").Value = -1120690.1

$ws = $wb.Worksheets.Item("BSM")
# Row 13
$ws.Range("H13").Value = 57450
$ws.Range("J13").Value = 57450
$ws.Range("L13").Value = 57450
$ws.Range("N13").Value = -57786
# Row 96
$ws.Range("H96").Value = 51142.668
$ws.Range("I96").Value = 21714
$ws.Range("K96").Value = 21714
$ws.Range("M96").Value = -18968
# Row 134
$ws.Range("H134").Value = 28544.357
$ws.Range("I134").Value = 1430.6451
$ws.Range("K134").Value = 4291.9353
$ws.Range("M134").Value = -1756.9353

$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 814.9167
$ws.Range("J94").Value = 1107.4286
$ws.Range("L94").Value = 1107.4286
$ws.Range("N94").Value = -2009.4286
# Row 133
$ws.Range("H133").Value = 55899.8
$ws.Range("J133").Value = 55899.8
$ws.Range("L133").Value = 55899.8
$ws.Range("N133").Value = -60959.8
# Row 134
$ws.Range("H134").Value = 598531.0600000001
$ws.Range("I134").Value = 4113
$ws.Range("K134").Value = 12339
$ws.Range("M134").Value = -9804

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 870.4
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 1075.5
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 3226.5
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -3450.5
# Row 34
$ws.Range("H34").Value = 6959.8857
$ws.Range("I34").Value = 621.3570999999999
$ws.Range("K34").Value = 1864.0713
$ws.Range("M34").Value = -1780.0713
# Row 68
$ws.Range("H68").Value = 1799.7213
$ws.Range("I68").Value = 1399.875
$ws.Range("K68").Value = 4199.625
$ws.Range("M68").Value = -3388.625
# Row 71
$ws.Range("H71").Value = 1799.7213
$ws.Range("I71").Value = 1399.875
$ws.Range("K71").Value = 12598.875
$ws.Range("M71").Value = -8542.875
# Row 107
$ws.Range("H107").Value = 23650.17
$ws.Range("I107").Value = 712.2857
$ws.Range("J107").Value = 27664.3
$ws.Range("K107").Value = 2136.8571
$ws.Range("L107").Value = 82992.89999999999
$ws.Range("M107").Value = -216.8571000000002
$ws.Range("N107").Value = -86832.89999999999
# Row 129
$ws.Range("H129").Value = 1958.2
$ws.Range("J129").Value = 2373
$ws.Range("L129").Value = 7119
$ws.Range("N129").Value = -17119
# Row 137
$ws.Range("H137").Value = 2102.5
$ws.Range("J137").Value = 3092
$ws.Range("L137").Value = 9276
$ws.Range("N137").Value = -19476
# Row 140
$ws.Range("H140").Value = 3382.7144
$ws.Range("I140").Value = 2169.4666
$ws.Range("K140").Value = 6508.399800000001
$ws.Range("M140").Value = -1328.399800000001

$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Range("H49").Value = 29929
$ws.Range("J49").Value = 29929
$ws.Range("L49").Value = 29929
$ws.Range("N49").Value = -30297

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 4000
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -3705
$ws.Range("N22").Value = -4590
# Row 27
$ws.Range("H27").Value = 4000
$ws.Range("I27").Value = 4000
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 4000
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = -3893
$ws.Range("N27").Value = -4214
# Row 32
$ws.Range("H32").Value = 22703.334
$ws.Range("I32").Value = 9055
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 9055
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -8738
$ws.Range("N32").Value = -50634
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = 0
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = 0
# Row 74
$ws.Range("H74").Value = 35239.4
$ws.Range("I74").Value = 30197
$ws.Range("K74").Value = 30197
$ws.Range("M74").Value = -29199
# Row 77
$ws.Range("H77").Value = 35239.4
$ws.Range("I77").Value = 30197
$ws.Range("K77").Value = 90591
$ws.Range("M77").Value = -85599
# Row 100
$ws.Range("H100").Value = 68702.94
$ws.Range("I100").Value = 107580
$ws.Range("K100").Value = 107580
$ws.Range("M100").Value = -107039
# Row 122
$ws.Range("H122").Value = 791492.0600000001
$ws.Range("I122").Value = 4647.5
$ws.Range("K122").Value = 13942.5
$ws.Range("M122").Value = -11492.5
# Row 132
$ws.Range("H132").Value = 6566.304
$ws.Range("I132").Value = 5220.3125
$ws.Range("K132").Value = 15660.9375
$ws.Range("M132").Value = -13130.9375

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 59831.94
$ws.Range("I107").Value = 101019.2
$ws.Range("J107").Value = 993
$ws.Range("K107").Value = 303057.6
$ws.Range("L107").Value = 2979
$ws.Range("M107").Value = -301137.6
$ws.Range("N107").Value = -6819
# Row 126
$ws.Range("H126").Value = 11814.091
$ws.Range("I126").Value = 2666.6667
$ws.Range("K126").Value = 8000.000100000001
$ws.Range("M126").Value = -5530.000100000001
# Row 132
$ws.Range("H132").Value = 27341.543
$ws.Range("I132").Value = 1935.5714
$ws.Range("J132").Value = 66861.94500000001
$ws.Range("K132").Value = 5806.7142
$ws.Range("L132").Value = 200585.835
$ws.Range("M132").Value = -3276.7142
$ws.Range("N132").Value = -205645.835
